$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "債務" (Debt) -- was sheet5.xml
# Fix the header row (previously duplicated the data row's own values)
# and append the English property_category/category/date/... columns
# (H:N) that every other sheet in the workbook already carries.
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")

# Header row 1 (B1:G1) -- correct field names instead of the
# mistakenly-duplicated data values.
$debt.Range("B1").Value = "species"
$debt.Range("C1").Value = "debtor"
$debt.Range("D1").Value = "owner"
$debt.Range("E1").Value = "total"
$debt.Range("F1").Value = "register_date"
$debt.Range("G1").Value = "register_reason"

# E2 becomes a real number instead of a text string.
$debt.Range("E2").Value = 5000000

# Extend header formatting (bold + border, matching B1:G1) across H1:N1
# before filling in the new header labels.
$debt.Range("G1").Copy() | Out-Null
$debt.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$debt.Range("H1").Value = "property_category"
$debt.Range("I1").Value = "category"
$debt.Range("J1").Value = "date"
$debt.Range("K1").Value = "legislator_name"
$debt.Range("L1").Value = "legislator_id"
$debt.Range("M1").Value = "source_file"
$debt.Range("N1").Value = "index"

# New data cells H2:N2. J2 ("date") must stay plain text ("2013-12-31")
# rather than be auto-converted to a date serial, so force a text
# number format on that cell first (matches the other sheets, which
# all store this same column as shared-string text).
$debt.Range("J2").NumberFormat = "@"

$debt.Range("H2").Value = "debt"
$debt.Range("I2").Value = "normal"
$debt.Range("J2").Value = "2013-12-31"
$debt.Range("K2").Value = "楊麗環"
$debt.Range("L2").Value = 960
$debt.Range("M2").Value = "tmp4fed1"
$debt.Range("N2").Value = 104

# ---------------------------------------------------------------------
# Sheet "事業投資" (Business investment) -- was sheet6.xml
# Same fix: correct header row, append H:N columns for every row.
# ---------------------------------------------------------------------
$inv = $wb.Worksheets.Item("事業投資")

$inv.Range("B1").Value = "owner"
$inv.Range("C1").Value = "company"
$inv.Range("D1").Value = "address"
$inv.Range("E1").Value = "total"
$inv.Range("F1").Value = "register_date"
$inv.Range("G1").Value = "register_reason"

# Extend header formatting across H1:N1 then fill in labels.
$inv.Range("G1").Copy() | Out-Null
$inv.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$inv.Range("H1").Value = "property_category"
$inv.Range("I1").Value = "category"
$inv.Range("J1").Value = "date"
$inv.Range("K1").Value = "legislator_name"
$inv.Range("L1").Value = "legislator_id"
$inv.Range("M1").Value = "source_file"
$inv.Range("N1").Value = "index"

# J2:J4 ("date") must stay plain text ("2013-12-31") instead of being
# auto-converted to date serials -- force text format first.
$inv.Range("J2:J4").NumberFormat = "@"

# Row 2 (index 109)
$inv.Range("H2").Value = "investment"
$inv.Range("I2").Value = "normal"
$inv.Range("J2").Value = "2013-12-31"
$inv.Range("K2").Value = "楊麗環"
$inv.Range("L2").Value = 960
$inv.Range("M2").Value = "tmp4fed1"
$inv.Range("N2").Value = 109

# Row 3 (index 110)
$inv.Range("E3").Value = 2000000
$inv.Range("H3").Value = "investment"
$inv.Range("I3").Value = "normal"
$inv.Range("J3").Value = "2013-12-31"
$inv.Range("K3").Value = "楊麗環"
$inv.Range("L3").Value = 960
$inv.Range("M3").Value = "tmp4fed1"
$inv.Range("N3").Value = 110

# Row 4 (index 112)
$inv.Range("E4").Value = 2500000
$inv.Range("H4").Value = "investment"
$inv.Range("I4").Value = "normal"
$inv.Range("J4").Value = "2013-12-31"
$inv.Range("K4").Value = "楊麗環"
$inv.Range("L4").Value = 960
$inv.Range("M4").Value = "tmp4fed1"
$inv.Range("N4").Value = 112
